# add Model 6 for ID 471
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update G8: week selection now covers "1 week" + "2 week" ---
$ws.Range("G8").Value = "1 week + `n2 week"
$ws.Range("G8").WrapText = $true

# --- Append new row 9: Model 6, Speed ID 471 ---
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 471
$ws.Range("C9").Value = "scaled speed weeakday o.h daypart o.h "
$ws.Range("D9").Value = "lstm(50)+do(.3)`nlstm/50)+do(.3)`nlstm/33)"
$ws.Range("E9").Value = 30
$ws.Range("F9").Value = "3h back`n0h forward"
$ws.Range("G9").Value = "1 week + `n2 week"
$ws.Range("H9").Value = "March`nApril`nMay"
$ws.Range("I9").Value = "First 7 days of June"
$ws.Range("J9").Value = 12.69
$ws.Range("K9").Value = 14.6
$ws.Range("L9").Value = 23.29
$ws.Range("M9").Value = "Adding a two prev week hasn" + [char]0x2019 + "t change loss"

$ws.Range("C9:D9").WrapText = $true
$ws.Range("F9:I9").WrapText = $true
$ws.Range("J9:K9").WrapText = $true
$ws.Range("M9").WrapText = $true

$ws.Rows.Item(9).RowHeight = 33

# --- Theme colors: swap dk1/lt1 (black/white -> yellow/black) ---
$theme = $wb.Theme
$cs = $theme.ThemeColorScheme
$cs.Colors(1).RGB = 65535   # dk1 -> FFFF00
$cs.Colors(2).RGB = 0       # lt1 -> 000000

# --- View: scroll / selection ---
$excel.Goto($ws.Range("B7"), $true)
$ws.Range("M9").Select()
